# HOT-56 - Docs ready to convert to HTML.
#
# 1. Insert two new paragraphs right after "Regular SQL selects shares..."
#    (before the "How to Run this example" Heading2 paragraph):
#      - an empty Normal paragraph
#      - a Normal paragraph with text about <column> tags, where "<column>"
#        uses the "Source" character style, and a bookmark
#        (__DdeLink__176_1849005264) wraps the first sentence.
# 2. The pre-existing bookmark __DdeLink__1133_7579575171 shifts from
#    id 0 -> id 1 automatically because the new bookmark takes id 0.
# 3/4. Update the footer's cached PAGE / NUMPAGES field results (1 -> 3,
#      2 -> 3) to reflect the now longer document.

$d = $word.ActiveDocument

# --- Locate the anchor paragraph ("Regular SQL selects shares...") ---
$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Regular SQL selects shares the syntax*") {
        $anchorIndex = $i
        break
    }
}

$anchor = $d.Paragraphs.Item($anchorIndex)

# --- Insert the first (empty) new paragraph ---
$anchor.Range.InsertParagraphAfter()

# --- Insert the second new paragraph (still empty at this point) ---
$emptyPara = $d.Paragraphs.Item($anchorIndex + 1)
$emptyPara.Range.InsertParagraphAfter()

# --- Fill the second new paragraph with text ---
$textPara = $d.Paragraphs.Item($anchorIndex + 2)
$paraStart = $textPara.Range.Start
$fullText = "SQL selects can include <column> tags to force the Java name and Java type of the result set. When included they need to be included outside the CDATA section, if any."
$textPara.Range.Text = $fullText

# --- Apply the "Source" character style to the "<column>" run ---
$colTag = "<column>"
$colStart = $paraStart + $fullText.IndexOf($colTag)
$colEnd = $colStart + $colTag.Length
$colRange = $d.Range($colStart, $colEnd)
$colRange.Style = "Source"

# --- Wrap the first sentence in the new bookmark ---
$bmEndMarker = "result set."
$bmEnd = $paraStart + $fullText.IndexOf($bmEndMarker) + $bmEndMarker.Length
$bmRange = $d.Range($paraStart, $bmEnd)
$d.Bookmarks.Add("__DdeLink__176_1849005264", $bmRange)

# --- Update the footer's cached PAGE / NUMPAGES field results ---
$footer = $d.Sections.Item(1).Footers.Item(1)
$fields = $footer.Range.Fields

$pageField = $fields.Item(1)
$pageResult = $pageField.Result
$pageResult.Find.Execute("1", $false, $false, $false, $false, $false, $true, 0, $false, "3", 1)

$footer2 = $d.Sections.Item(1).Footers.Item(1)
$fields2 = $footer2.Range.Fields
$numPagesField = $fields2.Item(2)
$numPagesResult = $numPagesField.Result
$numPagesResult.Find.Execute("2", $false, $false, $false, $false, $false, $true, 0, $false, "3", 1)
